# Slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") contains a 2-column table
# (the 3rd shape on the slide) whose table style is changed from the
# built-in "Medium Style 2 - Accent 1" ({EB24D34C-70CC-4A7E-83C3-786B10BD2F4B})
# to "Medium Style 2 - Accent 2" ({AB74CA05-EAB1-4CCD-933C-9F239ADE6B76}).
#
# This mirrors selecting the table (Table Tools > Design ribbon) and
# clicking the new style thumbnail in the Table Styles gallery, which PowerPoint
# implements as Table.ApplyStyle(StyleId, VisualStyleOnly).

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(16)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

$tbl.ApplyStyle("{AB74CA05-EAB1-4CCD-933C-9F239ADE6B76}", $true)
